$d = $word.ActiveDocument

# 1. Trim the old "To Do o a Done." ending off the Drag & Drop / SortableJS
#    bullet, leaving "... y añadirse a" (no trailing space) as its own run.
$found = $d.Content.Find.Execute(
    "añadirse a To Do o a Done.", $false, $false, $false, $false, $false,
    $true, 1, $false, "añadirse a", 2)

# 2. Locate that paragraph again so we can append a new run
#    " diferentes estados" right before its paragraph mark.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*arrastrarse*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range
$insertPoint = $d.Range($r.End - 1, $r.End - 1)
$insertPoint.InsertAfter(" diferentes estados")

# 3. Insert a new list-item paragraph right after the bullet (it inherits the
#    same pStyle/numPr as the Drag & Drop bullet), and fill in its text.
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($targetIndex + 1)
$p1.Range.InsertBefore("Nuevo estado IN PROGRESS añadido. Se elimina lógica antigua del campo Completada, y se añade nueva lógica creando campo Estado en BBDD.")

# 4. Insert a second new list-item paragraph right after that one.
$p1 = $d.Paragraphs.Item($targetIndex + 1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($targetIndex + 2)
$p2.Range.InsertBefore("Nueva funcionalidad editar tarea")
